{"js": "// Remove the trailing page-break paragraph and the empty paragraph that\n// followed it (the paragraphs right before the final sectPr), restoring\n// the document to end right after the \"JB22S953 \u8c37\u53e3\u58ee\u5e73\" author line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\n\n// Delete the last paragraph (the empty one) and the one before it\n// (the one holding the page-break run). Deleting from the end keeps\n// the remaining indices valid.\nif (count >= 2) {\n  paragraphs.items[count - 1].delete();\n  paragraphs.items[count - 2].delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing page-break paragraph (and the empty paragraph that\n# follows it, right before the final section properties) so the document\n# ends immediately after the \"JB22S953 \u8c37\u53e3\u58ee\u5e73\" author line.\n\n$d = $word.ActiveDocument\n\n# Locate the page break (^m / Chr(12)) in the document body.\n$r = $d.Content\n$found = $r.Find.Execute(\"^m\")\n\nif ($found) {\n    $breakStart = $r.Start\n    $breakEnd = $r.End\n\n    $count = $d.Paragraphs.Count\n    $idx = -1\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($breakStart -ge $p.Range.Start -and $breakStart -lt $p.Range.End) {\n            $idx = $i\n        }\n    }\n\n    if ($idx -gt 0) {\n        $startPos = $d.Paragraphs.Item($idx).Range.Start\n        $endIdx = $idx\n\n        # Also swallow the very next paragraph if it is empty (just a\n        # paragraph mark) \u2014 that is the stray blank paragraph left right\n        # before the sectPr.\n        if ($idx -lt $count) {\n            $nextPara = $d.Paragraphs.Item($idx + 1)\n            if ($nextPara.Range.Text -eq [char]13) {\n                $endIdx = $idx + 1\n            }\n        }\n\n        $endPos = $d.Paragraphs.Item($endIdx).Range.End\n        $d.Range($startPos, $endPos).Delete() | Out-Null\n    }\n}\n"}
